# Apply the "finish first version of pic" edit.
# Target sheet is "original task" (the active sheet, tab index 3).
#
# The two literal constants used throughout column L's formulas
# (15813 and 3671) are replaced with new values (2233 and 517).
# This also incidentally narrows the declared shared-formula range
# for si="4" from L10:L22 down to L10:L15 (it always only actually
# applied to L10:L15; L16:L22 use a separate formula/si="5").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("original task")
$ws.Activate()

# --- Column L block #1 (rows 2-8): 2.541579*15813/(3500/K#) -> 2.541579*2233/(3500/K#)
$ws.Range("L2").Formula  = "=2.541579*2233/(3500/K2)"
$ws.Range("L3").Formula  = "=2.541579*2233/(3500/K3)"
$ws.Range("L4").Formula  = "=2.541579*2233/(3500/K4)"
$ws.Range("L5").Formula  = "=2.541579*2233/(3500/K5)"
$ws.Range("L6").Formula  = "=2.541579*2233/(3500/K6)"
$ws.Range("L7").Formula  = "=2.541579*2233/(3500/K7)"
$ws.Range("L8").Formula  = "=2.541579*2233/(3500/K8)"

# --- Column L block #2 (rows 9-15): 3*1.4+15813*1.143/(3500/K#) -> 3*1.4+2233*1.143/(3500/K#)
$ws.Range("L9").Formula  = "=3*1.4+2233*1.143/(3500/K9)"
$ws.Range("L10").Formula = "=3*1.4+2233*1.143/(3500/K10)"
$ws.Range("L11").Formula = "=3*1.4+2233*1.143/(3500/K11)"
$ws.Range("L12").Formula = "=3*1.4+2233*1.143/(3500/K12)"
$ws.Range("L13").Formula = "=3*1.4+2233*1.143/(3500/K13)"
$ws.Range("L14").Formula = "=3*1.4+2233*1.143/(3500/K14)"
$ws.Range("L15").Formula = "=3*1.4+2233*1.143/(3500/K15)"

# --- Column L block #3 (rows 16-22): 3*1.4+(15813-3671)*1.143/(3500/K#) -> 3*1.4+(2233-517)*1.143/(3500/K#)
$ws.Range("L16").Formula = "=3*1.4+(2233-517)*1.143/(3500/K16)"
$ws.Range("L17").Formula = "=3*1.4+(2233-517)*1.143/(3500/K17)"
$ws.Range("L18").Formula = "=3*1.4+(2233-517)*1.143/(3500/K18)"
$ws.Range("L19").Formula = "=3*1.4+(2233-517)*1.143/(3500/K19)"
$ws.Range("L20").Formula = "=3*1.4+(2233-517)*1.143/(3500/K20)"
$ws.Range("L21").Formula = "=3*1.4+(2233-517)*1.143/(3500/K21)"
$ws.Range("L22").Formula = "=3*1.4+(2233-517)*1.143/(3500/K22)"

# --- Selection moves to L22 (matches the new <selection> in the diff)
$ws.Range("L22").Select()

$wb.Save()
